$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.537.45"
$ws.Range("E2").Value = "  +3.09%  "
$ws.Range("D3").Value = "2.548.56"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.59"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.83"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("D9").Value = "2.550.89"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.15"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "3.013.15"
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "68.488.64"
$ws.Range("E17").Value = "  +3.10%  "
$ws.Range("D18").Value = "2.538.16"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.05"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.61"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.39"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.76"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("E25").Value = "  +3.69%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.19"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("D28").Value = "2.679.70"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.11%  "
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "546.85"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.31"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.78%  "
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.29"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.97"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.69"
$ws.Range("D40").ClearFormats()
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.358"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.26"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("E44").Value = "  +2.30%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0285"
$ws.Range("E46").Value = "  +4.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.64"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.567"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.75"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  +1.21%  "
